# Rename the original sheet and add the new "IR" sheet after it.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Ultrasonic"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "IR"

# --- Header row ---
$ws2.Range("A1").Value = "Voltage"
$ws2.Range("B1").Value = "10 cm"
$ws2.Range("C1").Value = "20 cm"
$ws2.Range("D1").Value = "30 cm"
$ws2.Range("E1").Value = "40 cm"
$ws2.Range("F1").Value = "50 cm"
$ws2.Range("G1").Value = "60 cm"

# --- Measurement data (rows 2-21) ---
$irData = New-Object 'object[,]' 20,6
$irData[0,0] = 2.22
$irData[0,1] = 1.19
$irData[0,2] = 0.84
$irData[0,3] = 0.64
$irData[0,4] = 0.62
$irData[0,5] = 0.59
$irData[1,0] = 2.23
$irData[1,1] = 1.19
$irData[1,2] = 0.82
$irData[1,3] = 0.65
$irData[1,4] = 0.62
$irData[1,5] = 0.6
$irData[2,0] = 2.22
$irData[2,1] = 1.21
$irData[2,2] = 0.82
$irData[2,3] = 0.65
$irData[2,4] = 0.62
$irData[2,5] = 0.6
$irData[3,0] = 2.22
$irData[3,1] = 1.21
$irData[3,2] = 0.82
$irData[3,3] = 0.65
$irData[3,4] = 0.58
$irData[3,5] = 0.65
$irData[4,0] = 2.22
$irData[4,1] = 1.21
$irData[4,2] = 0.82
$irData[4,3] = 0.65
$irData[4,4] = 0.62
$irData[4,5] = 0.57
$irData[5,0] = 2.23
$irData[5,1] = 1.2
$irData[5,2] = 0.82
$irData[5,3] = 0.65
$irData[5,4] = 0.65
$irData[5,5] = 0.56
$irData[6,0] = 2.22
$irData[6,1] = 1.21
$irData[6,2] = 0.82
$irData[6,3] = 0.65
$irData[6,4] = 0.58
$irData[6,5] = 0.6
$irData[7,0] = 2.22
$irData[7,1] = 1.21
$irData[7,2] = 0.82
$irData[7,3] = 0.65
$irData[7,4] = 0.58
$irData[7,5] = 0.62
$irData[8,0] = 2.22
$irData[8,1] = 1.21
$irData[8,2] = 0.82
$irData[8,3] = 0.66
$irData[8,4] = 0.6
$irData[8,5] = 0.6
$irData[9,0] = 2.23
$irData[9,1] = 1.21
$irData[9,2] = 0.82
$irData[9,3] = 0.65
$irData[9,4] = 0.58
$irData[9,5] = 0.62
$irData[10,0] = 2.23
$irData[10,1] = 1.21
$irData[10,2] = 0.82
$irData[10,3] = 0.65
$irData[10,4] = 0.58
$irData[10,5] = 0.6
$irData[11,0] = 2.22
$irData[11,1] = 1.21
$irData[11,2] = 0.82
$irData[11,3] = 0.63
$irData[11,4] = 0.58
$irData[11,5] = 0.6
$irData[12,0] = 2.23
$irData[12,1] = 1.21
$irData[12,2] = 0.82
$irData[12,3] = 0.65
$irData[12,4] = 0.58
$irData[12,5] = 0.6
$irData[13,0] = 2.23
$irData[13,1] = 1.21
$irData[13,2] = 0.82
$irData[13,3] = 0.65
$irData[13,4] = 0.57
$irData[13,5] = 0.6
$irData[14,0] = 2.23
$irData[14,1] = 1.21
$irData[14,2] = 0.82
$irData[14,3] = 0.65
$irData[14,4] = 0.57
$irData[14,5] = 0.62
$irData[15,0] = 2.23
$irData[15,1] = 1.21
$irData[15,2] = 0.82
$irData[15,3] = 0.65
$irData[15,4] = 0.58
$irData[15,5] = 0.6
$irData[16,0] = 2.23
$irData[16,1] = 1.21
$irData[16,2] = 0.82
$irData[16,3] = 0.66
$irData[16,4] = 0.58
$irData[16,5] = 0.6
$irData[17,0] = 2.23
$irData[17,1] = 1.21
$irData[17,2] = 0.82
$irData[17,3] = 0.66
$irData[17,4] = 0.58
$irData[17,5] = 0.62
$irData[18,0] = 2.23
$irData[18,1] = 1.21
$irData[18,2] = 0.83
$irData[18,3] = 0.65
$irData[18,4] = 0.58
$irData[18,5] = 0.62
$irData[19,0] = 2.23
$irData[19,1] = 1.21
$irData[19,2] = 0.82
$irData[19,3] = 0.65
$irData[19,4] = 0.58
$irData[19,5] = 0.6
$ws2.Range("B2:G21").Value = $irData

# --- Averages / Standard Dev rows ---
$ws2.Range("A23").Value = "Averages"
$ws2.Range("B23").Formula = "=AVERAGE(B2:B21)"
$ws2.Range("C23:G23").Formula = "=AVERAGE(C2:C21)"

$ws2.Range("A24").Value = "Standard Dev"
$ws2.Range("B24").Formula = "=STDEV(B2:B21)"
$ws2.Range("C24:G24").Formula = "=STDEV(C2:C21)"

# --- Selections / active sheet ---
$ws1.Range("I6").Select()
$ws2.Range("P17").Select()

